$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")
$ws2 = $wb.Worksheets.Item("Лаба1")

$ws2.Visible = $true

$ws1.Range("B4").Value = 5
$ws1.Range("B10").Value = 5
$ws1.Range("B11").Value = "+"
$ws1.Range("B14").Value = "*"
$ws1.Range("B18").Value = 5
$ws1.Range("B20").Value = 5

$ws1.Range("B21").Select()

$ws2.Activate()
$ws2.Range("B8").Select()

$ws1.Activate()
